$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 20 de Abril de 2020 a las 08:22"

# Swap the country names for rows 123 and 124 (Paraguay <-> El Salvador)
$ws.Range("A123").Value = "El Salvador"
$ws.Range("A124").Value = "Paraguay"

# Update statistic values per country row (Noruega)
$ws.Range("B36").Value = 7103
$ws.Range("C36").Value = 25
$ws.Range("E36").Value = 6906

# Tailandia
$ws.Range("B57").Value = 2792
$ws.Range("C57").Value = 27
$ws.Range("D57").Value = 1999
$ws.Range("E57").Value = 746

# Hungria
$ws.Range("B62").Value = 1984
$ws.Range("C62").Value = 68
$ws.Range("D62").Value = 267
$ws.Range("E62").Value = 1518
$ws.Range("F62").Value = 60
$ws.Range("G62").Value = 10
$ws.Range("H62").Value = 199

# Bulgaria
$ws.Range("B85").Value = 915
$ws.Range("C85").Value = 21
$ws.Range("D85").Value = 167
$ws.Range("E85").Value = 705
$ws.Range("G85").Value = 1
$ws.Range("H85").Value = 43

# Taiwan
$ws.Range("B107").Value = 422
$ws.Range("C107").Value = 2
$ws.Range("D107").Value = 203
$ws.Range("E107").Value = 213

# Row 123 (now El Salvador)
$ws.Range("B123").Value = 218
$ws.Range("C123").Value = 17
$ws.Range("D123").Value = 46
$ws.Range("E123").Value = 165
$ws.Range("F123").Value = 2
$ws.Range("H123").Value = 7

# Row 124 (now Paraguay)
$ws.Range("B124").Value = 208
$ws.Range("C124").Value = 2
$ws.Range("D124").Value = 46
$ws.Range("E124").Value = 154
$ws.Range("H124").Value = 8
